# Add a new "Interact" column (I) to Sheet1: an English key row, a
# Chinese localized label row, and a 0/1 flag for each data row telling
# whether that object/unit can be interacted with.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the formatting from column H (same header/body style already used
# across the table) onto the new column I before filling in values.
$ws.Range("H1:H17").Copy()
$ws.Range("I1:I17").PasteSpecial(-4122) # xlPasteFormats

# Header rows (row 1 = English key, row 2 = Chinese localized label).
$ws.Range("I1").Value = "Interact"
$ws.Range("I2").Value = "交互"

# Data rows 3-17: whether the unit/object can be interacted with.
$values = @(0, 0, 0, 1, 1, 1, 1, 1, 1, 1, 1, 1, 0, 0, 0)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 9).Value = $values[$i]
}

# Update the selection to match the editor's last cursor position.
$ws.Range("I20").Select()
